# "update plots for each sample"
# The marker CYP2D6_002 / CYP2D6_10B wild-type allele for sample S1 was
# previously undetected (peak height 1000 below the min_height threshold of
# 1000). Re-running the detection with an updated min_height of 800 now
# finds the wild-type peak, flipping the genotype call for this sample from
# homozygous-mutant (TT) to heterozygous (CT), and the overall genotype
# result from *10B/*10B to *1/*10B.

$wb = $excel.ActiveWorkbook

# --- peak_table (sheet1): lower the wild-type peak-height threshold for
#     the CYP2D6_10B marker (row 3: CYP2D6_002 / CYP2D6_10B) from 1000 to 800
$ws1 = $wb.Worksheets.Item("peak_table")
$ws1.Cells.Item(3, 14).Value = 800   # N3 = w_height

# --- allele_table (sheet2): the wild-type allele row (row 4) now detects
#     a peak given the lower min_height threshold
$ws2 = $wb.Worksheets.Item("allele_table")
$ws2.Cells.Item(4, 11).Value = 800        # K4 = min_height
$ws2.Cells.Item(4, 13).Value = $true       # M4 = is_detected
$ws2.Cells.Item(4, 14).Value = 35         # N4 = peak
$ws2.Cells.Item(4, 15).Value = 31.88       # O4 = size
$ws2.Cells.Item(4, 16).Value = 847         # P4 = height
$ws2.Cells.Item(4, 17).Value = "ok"        # Q4 = status
$ws2.Cells.Item(4, 18).Value = ""          # R4 = message (cleared)

# --- marker_table (sheet3): genotype call for this marker changes from
#     homozygous mutant (TT) to heterozygous (CT)
$ws3 = $wb.Worksheets.Item("marker_table")
$ws3.Cells.Item(3, 7).Value = "CT"          # G3 = genotype
$ws3.Cells.Item(3, 8).Value = "heterozygous" # H3 = phenotype

# --- genotype_result (sheet4): overall sample genotype updates accordingly
$ws4 = $wb.Worksheets.Item("genotype_result")
$ws4.Cells.Item(2, 2).Value = "*1/*10B"     # B2 = genotype
